# Weekly update: insert this week's new record as a new row right after
# the header (row 64), shifting all the existing historical rows (64-149)
# down by one (they become rows 65-150). This mirrors Excel's default
# "Insert" behaviour on EntireRow (shift down).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 64; everything below (old 64..149) shifts
# down to 65..150, and the sheet's used range grows from R149 to R150.
$ws.Rows(64).Insert()

# Populate the newly inserted row 64 with this week's data.
$ws.Range("A64").Value = 4
$ws.Range("B64").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C64").Value = "Los Lagos"
$ws.Range("D64").Value = 44895
$ws.Range("E64").Value = 10
$ws.Range("F64").Value = 100112052
$ws.Range("G64").Value = "Albahaca"
$ws.Range("H64").Value = "Sin especificar"
$ws.Range("I64").Value = "Primera"
$ws.Range("J64").Value = 60
$ws.Range("K64").Value = 8000
$ws.Range("L64").Value = 8000
$ws.Range("M64").Value = 8000
$ws.Range("N64").Value = "$/docena de matas"
$ws.Range("O64").Value = "Región Metropolitana"
$ws.Range("P64").Value = 1333
$ws.Range("Q64").Value = 6
$ws.Range("R64").Value = "Hortaliza"
